$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two brand-new columns at the very front of the table:
#    "Factory" (A) and "Buy Plan Season" (B). Everything that used to live
#    in columns A:R shifts two columns to the right (becomes C:T).
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.Insert()

$ws.Range("A1").Value = "Factory"
$ws.Range("B1").Value = "Buy Plan Season"
$ws.Range("A2").Value = "&=result.Factory"
$ws.Range("B2").Value = "&=result.BuyPlanSeason"

# Match the look of the header / data rows to the (shifted) neighbour column
# so the new columns render identically to the rest of the header/data rows.
$ws.Range("A1:B1").Font.Size = $ws.Range("C1").Font.Size()
$ws.Range("A1:B1").Font.Name = $ws.Range("C1").Font.Name()
$ws.Range("A1:B1").Interior.Color = $ws.Range("C1").Interior.Color()
$ws.Range("A1:B1").HorizontalAlignment = $ws.Range("C1").HorizontalAlignment()
$ws.Range("A1:B1").VerticalAlignment = $ws.Range("C1").VerticalAlignment()

$ws.Range("A2:B2").Font.Size = $ws.Range("C2").Font.Size()
$ws.Range("A2:B2").Font.Name = $ws.Range("C2").Font.Name()

# Column widths for the two new columns (raw OOXML widths 15.42578125 /
# 22 expressed as the character-based ColumnWidth this host expects).
$ws.Columns("A").ColumnWidth = 14.666666666666666
$ws.Columns("B").ColumnWidth = 21.166666666666668

# ---------------------------------------------------------------------------
# 2. Move the "Remark" column (now at column O after the insert above) so it
#    becomes the very last column of the table instead of sitting in the
#    middle.
# ---------------------------------------------------------------------------
$ws.Columns("O").Cut() | Out-Null
$ws.Columns("U").Insert()

# ---------------------------------------------------------------------------
# 3. Refresh the autofilter so it spans the new full width of the table
#    (A1:T1 instead of A1:R1).
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:T1").AutoFilter() | Out-Null

# ---------------------------------------------------------------------------
# 4. Update the workbook-level _FilterDatabase defined name to match the new
#    autofilter range.
# ---------------------------------------------------------------------------
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$T`$1"

# ---------------------------------------------------------------------------
# 5. Restore the originally selected cell on the sheet.
# ---------------------------------------------------------------------------
$ws.Range("G14").Select() | Out-Null

Write-Output "done"
